$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '55.811.78'
$ws.Range('D3').Value = '2.970.15'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '502.11'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '136.36'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '3.485.09'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D14').Value = '25.83'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').Value = '55.879.03'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '2.972.52'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '5.96'
$ws.Range('D19').Value = '12.82'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').Value = '7.93'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').Value = '326.66'
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '64.27'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').Value = '3.094.43'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').Value = '0.162'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').Value = '0.0₃0886'
$ws.Range('E28').Value = '  -1.22%  '
$ws.Range('D29').Value = '6.33'
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('D30').Value = '6.94'
$ws.Range('E30').Value = '  -1.86%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('D34').Value = '153.25'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').Value = '4.46'
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('E36').Value = '  -1.91%  '
$ws.Range('D37').Value = '25.35'
$ws.Range('E37').Value = '  +4.76%  '
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').Value = '0.0655'
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('D40').Value = '3.009.92'
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').Value = '36.71'
$ws.Range('E41').Value = '  -3.08%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '3.75'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '0.648'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').Value = '2.156.38'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('D47').Value = '5.79'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').Value = '0.915'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').Value = '19.45'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('E51').Value = '  -3.65%  '
